$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "Exploratory Data Analysis of the ToothGrowth Dataset"
# Heading2 paragraph. A new "Assumptions" section must be inserted
# immediately before it.
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "Exploratory Data Analysis of the ToothGrowth Dataset") {
        $target = $i
        break
    }
}

# Remember the list template used by the existing "Compact" numbered
# bullets (the instructions list under Overview) so the new bullet
# list can be created as a sibling list (same abstract numbering,
# but its own fresh numId that restarts at 1).
$srcListPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.ListFormat.ListType -ne 0) {
        $srcListPara = $cand
        break
    }
}
$listTemplate = $srcListPara.Range.ListFormat.ListTemplate

# ------------------------------------------------------------------
# Insert the new paragraphs' text before the heading. InsertBefore on
# a collapsed range inserts new paragraphs that inherit the style of
# the paragraph it is anchored to (Heading2), which conveniently
# matches the first new paragraph ("Assumptions").
# ------------------------------------------------------------------
$r = $d.Paragraphs($target).Range
$r.Collapse(1)

$newText = "Assumptions`r" + `
    "This analysis and its conclusions are dependent on the following set of assumptions:`r" + `
    "The sample populations are ramdomly selected and independent of each other.`r" + `
    "The population was comprised of similar guinea pigs as was the methods for delivery of supplements, thus the variables are independent and identically distributed (iid).`r" + `
    "Variances of tooth growth are different when using different supplements and delivery methods.`r" + `
    "Tooth growth follows a normal distribution.`r" + `
    "A confidence interval of 95% is satisfactory for our conclusions.`r"

$r.InsertBefore($newText)

# Paragraph indices after insertion:
#   target      -> "Assumptions"                       (Heading2)
#   target + 1  -> "This analysis ..."                 (Normal)
#   target + 2  -> "The sample populations ..."         (Compact, bullet 1)
#   target + 3  -> "The population was comprised ..."   (Compact, bullet 2)
#   target + 4  -> "Variances of tooth growth ..."      (Compact, bullet 3)
#   target + 5  -> "Tooth growth follows ..."           (Compact, bullet 4)
#   target + 6  -> "A confidence interval of 95% ..."   (Compact, bullet 5)
#   target + 7  -> "Exploratory Data Analysis ..."      (Heading2, original)

$pAssumptions = $target
$pIntro       = $target + 1
$pBullet1     = $target + 2
$pBullet2     = $target + 3
$pBullet3     = $target + 4
$pBullet4     = $target + 5
$pBullet5     = $target + 6
$pHeading     = $target + 7

$d.Paragraphs($pIntro).Style = "Normal"
$d.Paragraphs($pBullet1).Style = "Compact"
$d.Paragraphs($pBullet2).Style = "Compact"
$d.Paragraphs($pBullet3).Style = "Compact"
$d.Paragraphs($pBullet4).Style = "Compact"
$d.Paragraphs($pBullet5).Style = "Compact"

# Apply the bullet numbering to the 5 new assumption bullets as a
# single contiguous range so they all share one freshly minted numId
# (the runtime mints a new w:num entry, cloned from the same abstract
# numbering definition, each time ApplyListTemplate targets content
# that isn't already tied to that exact num instance).
$bulletsRange = $d.Range($d.Paragraphs($pBullet1).Range.Start, $d.Paragraphs($pBullet5).Range.End)
$bulletsRange.ListFormat.ApplyListTemplate($listTemplate)

# ------------------------------------------------------------------
# Bookmarks: pandoc-style heading anchors are zero-width, sitting at
# the very start of the heading paragraph (bookmarkStart immediately
# followed by bookmarkEnd, before the run).
# ------------------------------------------------------------------
$assumptionsStart = $d.Paragraphs($pAssumptions).Range.Start
$bmRange = $d.Range($assumptionsStart, $assumptionsStart)
$d.Bookmarks.Add("assumptions", $bmRange)

Write-Host "Inserted Assumptions section before paragraph" $pHeading
